$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6
$ws.Range("A6").Value = 111742151
$ws.Range("B6").Value = 95524
$ws.Range("D6").Value = 'LC'
$ws.Range("E6").Value = 221944
$ws.Range("F6").Value = 'Lopplummer'
$ws.Range("G6").Value = 'Huperzia selago'
$ws.Range("H6").Value = '(L.) Bernh. ex Schrank & Mart.'
$ws.Range("Q6").Value = 331814.6184995985
$ws.Range("R6").Value = 6626778.67820756
$ws.Range("AC6").ClearContents()

# Row 7
$ws.Range("A7").Value = 111742077
$ws.Range("B7").Value = 78605
$ws.Range("D7").Value = 'LC'
$ws.Range("E7").Value = 6462
$ws.Range("F7").Value = 'Stuplav'
$ws.Range("G7").Value = 'Nephroma bellum'
$ws.Range("H7").Value = '(Spreng.) Tuck.'
$ws.Range("Q7").Value = 331735.1116598135
$ws.Range("R7").Value = 6626820.629936518
$ws.Range("AC7").Value = 'På rönn'

# Row 8
$ws.Range("A8").Value = 111742170
$ws.Range("B8").Value = 89369
$ws.Range("D8").Value = 'LC'
$ws.Range("E8").Value = 5447
$ws.Range("F8").Value = 'Vedticka'
$ws.Range("G8").Value = 'Fuscoporia viticola'
$ws.Range("H8").Value = '(Schwein.) Murrill'
$ws.Range("Q8").Value = 331846.7251686137
$ws.Range("R8").Value = 6626784.294692003
$ws.Range("AC8").ClearContents()

# Row 9
$ws.Range("A9").Value = 111742184
$ws.Range("B9").Value = 93159
$ws.Range("D9").Value = 'LC'
$ws.Range("E9").Value = 2666
$ws.Range("F9").Value = 'Grov fjädermossa'
$ws.Range("G9").Value = 'Neckera crispa'
$ws.Range("H9").Value = 'Hedw.'
$ws.Range("Q9").Value = 331833.6062344447
$ws.Range("R9").Value = 6626784.887086328
$ws.Range("AC9").Value = 'I bergsbrant'

# Row 10
$ws.Range("A10").Value = 111742070
$ws.Range("B10").Value = 78578
$ws.Range("D10").Value = 'NT'
$ws.Range("E10").Value = 6458
$ws.Range("F10").Value = 'Lunglav'
$ws.Range("G10").Value = 'Lobaria pulmonaria'
$ws.Range("H10").Value = '(L.) Hoffm.'
$ws.Range("Q10").Value = 331735.1116598135
$ws.Range("R10").Value = 6626820.629936518
$ws.Range("AC10").Value = 'På rönn'

# Row 11
$ws.Range("A11").Value = 111742101
$ws.Range("B11").Value = 94134
$ws.Range("D11").Value = 'NT'
$ws.Range("E11").Value = 53
$ws.Range("F11").Value = 'Vedtrappmossa'
$ws.Range("G11").Value = 'Crossocalyx hellerianus'
$ws.Range("H11").Value = '(Nees ex Lindenb.) Meyl.'
$ws.Range("Q11").Value = 331779.6127968954
$ws.Range("R11").Value = 6626798.429951042
$ws.Range("AC11").ClearContents()

# Row 12
$ws.Range("A12").Value = 111742181
$ws.Range("B12").Value = 93158
$ws.Range("D12").Value = 'LC'
$ws.Range("E12").Value = 2667
$ws.Range("F12").Value = 'Platt fjädermossa'
$ws.Range("G12").Value = 'Neckera complanata'
$ws.Range("H12").Value = '(Hedw.) Huebener'
$ws.Range("Q12").Value = 331833.6062344447
$ws.Range("R12").Value = 6626784.887086328
$ws.Range("AC12").Value = 'I bergsbrant'

# Row 13
$ws.Range("A13").Value = 111742096
$ws.Range("B13").Value = 94125
$ws.Range("D13").Value = 'NT'
$ws.Range("E13").Value = 54
$ws.Range("F13").Value = 'Skogstrappmossa'
$ws.Range("G13").Value = 'Anastrophyllum michauxii'
$ws.Range("H13").Value = '(F.Weber.) H.Buch'
$ws.Range("Q13").Value = 331779.6127968954
$ws.Range("R13").Value = 6626798.429951042
$ws.Range("AC13").Value = 'På både ved och på lodyta'
